$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header label changed: "Inject/Withdraw Rates and Costs" -> "Inject/Withdraw/Inventory Constraints"
# (the old string becomes unreferenced and is dropped from the shared-string table,
# which also re-numbers every other shared-string reference on the sheet automatically)
$ws.Range("K2").Value = "Inject/Withdraw/Inventory Constraints"

# Row 5
$ws.Range("M5").Value = 1356
$ws.Range("N5").Value = 1525

# Row 6
$ws.Range("M6").Value = 1306
$ws.Range("N6").Value = 1689

# Row 7
$ws.Range("M7").Value = 1206
$ws.Range("N7").Value = 1711

# Row 8
$ws.Range("M8").Value = 1145
$ws.Range("N8").Value = 1784

# Row 9
$ws.Range("M9").Value = 1005
$ws.Range("N9").Value = 1587

# Row 10
$ws.Range("L10").Value = 200
$ws.Range("M10").Value = 958
$ws.Range("N10").Value = 1604

# Row 11
$ws.Range("L11").Value = 400
$ws.Range("M11").Value = 910
$ws.Range("N11").Value = 1658

# Row 12
$ws.Range("K12").Value = 43709
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 1346
$ws.Range("N12").Value = 1479.25

# Row 13
$ws.Range("L13").Value = 300
$ws.Range("M13").Value = 1296
$ws.Range("N13").Value = 1638.33

# Row 14
$ws.Range("L14").Value = 500
$ws.Range("M14").Value = 1196
$ws.Range("N14").Value = 1659.6699999999998

# Row 15
$ws.Range("L15").Value = 600
$ws.Range("M15").Value = 1135
$ws.Range("N15").Value = 1730.48

# Row 16 - clear out the now-unused last data-table row
$ws.Range("K16").ClearContents()
$ws.Range("L16").ClearContents()

# Update the active selection to reflect where the user clicked last
$ws.Range("K3").Select()
